# Commit: "update h2 inputs with new SMR designs"
# Rename the header in A1 from "Reactor" to "Type" on the (single, active)
# worksheet, then leave the selection where the author's Excel session
# last left it (G20) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Type"

$ws.Range("G20").Select() | Out-Null
